$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in missing region values (Maine / Colorado / Massachusetts) for
# section-header rows, and the new site code "WP" for Wilcox Pond.
$ws.Range("B5").Value = "Maine"
$ws.Range("B8").Value = "Colorado"
$ws.Range("C10").Value = "WP"
$ws.Range("B18").Value = "Massachusetts"

# Update the active selection to reflect the last-edited cell.
$ws.Range("E33").Select()
